# fix courts upload logic
#
# The "latitude" column (B) for three courts (rows 4, 5, 7) was uploaded as
# raw numbers, which Excel silently reinterpreted as date serials (e.g.
# 7.8804 -> 2521818, formatted as a "mmm-yy" date) instead of keeping them
# as the decimal latitude strings used everywhere else in the column.
# Re-enter the values as text so they display/store correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: keep the existing bordered cell style, just switch it from a
#     date number format to plain text and re-enter the correct value.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "7.8804"

# --- Row 5: same data-entry bug, but this one also loses its cell border
#     and text wrapping in the fix (matches the new "plain" latitude style).
$ws.Range("B5").Borders.LineStyle = 0
$ws.Range("B5").WrapText = $false
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "12.9236"

# --- Row 7: same fix as row 5.
$ws.Range("B7").Borders.LineStyle = 0
$ws.Range("B7").WrapText = $false
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "12.5707"

# Leave the selection where the edit finished, on B7.
$ws.Range("B7").Select()
